$d = $word.ActiveDocument

# 1. Locate the paragraph that currently ends the "extensions" bullet list item
#    we want to add a new bullet after ("Extend client implementation to
#    support Android app").
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Extend client implementation to support Android app", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng = $find.Parent
$para = $rng.Paragraphs(1)

# 2. Split off a brand-new paragraph right after it, inheriting the same
#    paragraph/run formatting (ListParagraph style + numbering).
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newPara = $para.Next()
$newRng = $newPara.Range

# 3. Fill the new paragraph with the requested text. A sentinel character is
#    appended temporarily so that the later bookmark placement is not at the
#    very end of the paragraph (which the engine mis-serialises); it is
#    stripped again right before saving.
$newRng.InsertBefore("Extend client interaction with the device to allow multiple ways to introduce the destinationX")

# 4. Move the hidden "_GoBack" bookmark (tracks the last edit location) from
#    wherever it currently sits onto the newly typed text.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("destinationX", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markRng = $find2.Parent
$markRng.MoveStart(1, 11)

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $markRng)

# 5. Remove the sentinel character again, leaving the bookmark collapsed
#    right after "destination".
$markRng.Text = ""
